$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("ENCARGADO") assignments for each subtask -------------------
# Order matters: it controls the order new shared strings are interned in,
# which must match the target sharedStrings.xml (index 40 = "Jose Manuel",
# index 41 = "Enrique, Mario"). Setting D11 first, then D3 (before the plain
# "Enrique" cells) reproduces that ordering; the "Enrique" cells below all
# reuse the pre-existing shared string (index 23).

$ws.Range("D11").Value = "José Manuel"
$ws.Range("D3").Value = "Enrique, Mario"

$ws.Range("D4").Value = "Enrique"
$ws.Range("D6").Value = "Enrique"
$ws.Range("D8").Value = "Enrique"
$ws.Range("D9").Value = "Enrique"
$ws.Range("D10").Value = "Enrique"
$ws.Range("D13").Value = "Enrique"
$ws.Range("D14").Value = "Enrique"
$ws.Range("D19").Value = "Enrique"
$ws.Range("D20").Value = "Enrique"
$ws.Range("D24").Value = "Enrique, Mario"
$ws.Range("D26").Value = "Enrique"

# D13 and D24 pick up an explicit (automatic/black) font colour, which is
# what produces the extra font + cellXfs entry seen in the diff.
$ws.Range("D13").Font.Color = 0
$ws.Range("D24").Font.Color = 0

# These two summary rows (task-level "ENCARGADO" rollups) lose their values
# now that every subtask row carries its own owner.
$ws.Range("D12").Clear()
$ws.Range("D18").Clear()

# Restore the on-open selection/scroll state captured in the diff.
$ws.Range("E13").Select()
